$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# AMSIN sheet (sheet1): refine the precision of the existing B70
# timestamp, then append a new test-run row (row 71).
# ------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

$wsAmsin.Cells.Item(70, 2).Value = 45089.72860135417

$wsAmsin.Cells.Item(71, 1).Value = "'2023-07-31"
$bCell = $wsAmsin.Cells.Item(71, 2)
$bCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$bCell.Value = 45138.41657248842
$wsAmsin.Cells.Item(71, 3).Value = "180fnlrun"
$wsAmsin.Cells.Item(71, 4).Value = 105
$wsAmsin.Cells.Item(71, 5).Value = 104
$wsAmsin.Cells.Item(71, 6).Value = 1
$wsAmsin.Cells.Item(71, 7).Value = 3.32

# ------------------------------------------------------------------
# BETA sheet (sheet2): append two new test-run rows (34 and 35).
# ------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

$wsBeta.Cells.Item(34, 1).Value = "'2023-08-01"
$bCell = $wsBeta.Cells.Item(34, 2)
$bCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$bCell.Value = 45139.5682859375
$wsBeta.Cells.Item(34, 3).Value = "180beta"
$wsBeta.Cells.Item(34, 4).Value = 105
$wsBeta.Cells.Item(34, 5).Value = 64
$wsBeta.Cells.Item(34, 6).Value = 41
$wsBeta.Cells.Item(34, 7).Value = 11.58

$wsBeta.Cells.Item(35, 1).Value = "'2023-08-01"
$bCell = $wsBeta.Cells.Item(35, 2)
$bCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$bCell.Value = 45139.68982625
$wsBeta.Cells.Item(35, 3).Value = "180beta"
$wsBeta.Cells.Item(35, 4).Value = 105
$wsBeta.Cells.Item(35, 5).Value = 104
$wsBeta.Cells.Item(35, 6).Value = 1
$wsBeta.Cells.Item(35, 7).Value = 3.25

# ------------------------------------------------------------------
# AMS sheet (sheet3): append one new test-run row (36).
# ------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

$wsAms.Cells.Item(36, 1).Value = "'2023-08-01"
$bCell = $wsAms.Cells.Item(36, 2)
$bCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$bCell.Value = 45139.84663266376
$wsAms.Cells.Item(36, 3).Value = "180live"
$wsAms.Cells.Item(36, 4).Value = 105
$wsAms.Cells.Item(36, 5).Value = 105
$wsAms.Cells.Item(36, 6).Value = 0
$wsAms.Cells.Item(36, 7).Value = 2.86
